# EPBDS-6561 Add cycled dependency on beans
# Inserts a new "Address adr" row (a field of type Address named adr)
# right after the existing Person datatype properties (rows 2-6),
# pushing the Spreadsheet/Step/Environment block down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting rows 7+ down by one.
$ws.Rows("7").Insert()

# Populate the newly inserted row with the new datatype field.
$ws.Range("B7").Value = "Address"
$ws.Range("C7").Value = "adr"

# Match the resulting selection left behind in the saved file.
$ws.Range("C8").Select()
